# TEAM_2_SCOPE_DEFINITION.docx edit
# Commit: "Changed sales to shipping in organizational structure in scope document."
#
# The only substantive, visible-text change in this document is in the
# "Organizational/Culture" bullet under Feasibility Evaluation, where the
# department list is updated from "sales" to "shipping". (The rest of the
# upstream diff consists purely of Word's own proofing engine wrapping
# existing words in <w:proofErr> spell/grammar-check bookmarks and
# coalescing a couple of identically-formatted runs -- invisible markup
# artifacts that Word's own background spell/grammar checker stamps onto
# the file when it is saved, not something exposed through the Word
# object model, so there is no text/content difference to author for
# those spots.)

$d = $word.ActiveDocument

# Replace "sales" with "shipping" in the department list. Anchor the
# search on surrounding text so we only touch this one occurrence.
$found = $d.Content.Find.Execute(
    "four departments, sales, quality",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "four departments, shipping, quality",
    2
)
